# Mirrors the upstream "Moving from Apache POI 4.1.0 to 5.2.3" fix for
# tests/.../noBookmark/noBookmark-expected-generation.docx.
#
# The underlying OOXML diff between the two POI versions only touches
# incidental/serialization-level details (the literal ST_OnOff spelling
# used for a boolean-true <w:b> element, the run's w:rsidR save-session
# token, and the element order inside <w:rPr>) - none of which change the
# actual rendered content: both runs were bold before the upgrade and
# remain bold after it, and the REF-field result run keeps its
# <w:noProof/> marker. There is no text/content change to make.
#
# Re-assert the (already-true) Bold formatting on the two runs the diff
# touches, which is the content-level echo of that upgrade reachable
# through the Word object model.

$d = $word.ActiveDocument

# Run 1: "dangling reference for bookmark bookmark1" (bold, red)
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "dangling reference for bookmark bookmark1"
if ($find1.Execute()) {
    $find1.Parent.Font.Bold = $true
}

# Run 2: "a reference to bookmark1" (the REF field's cached result, bold + noProof)
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "a reference to bookmark1"
if ($find2.Execute()) {
    $find2.Parent.Font.Bold = $true
}
